$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 2437.375
$ws.Range("I8").Value = 1933.0667
$ws.Range("K8").Value = 5799.2001
$ws.Range("M8").Value = -5660.2001
$ws.Range("H29").Value = 7999
$ws.Range("J29").Value = 7999
$ws.Range("L29").Value = 23997
$ws.Range("N29").Value = -24559
$ws.Range("H32").Value = 4561.5386
$ws.Range("I32").Value = 4050
$ws.Range("K32").Value = 4050
$ws.Range("M32").Value = -3724
$ws.Range("H38").Value = 614.9
$ws.Range("I38").Value = 127.666664
$ws.Range("K38").Value = 382.999992
$ws.Range("M38").Value = -10.99999200000002
$ws.Range("H43").Value = 8140.4346
$ws.Range("I43").Value = 8265.6
$ws.Range("J43").Value = 7905.75
$ws.Range("K43").Value = 8265.6
$ws.Range("L43").Value = 7905.75
$ws.Range("M43").Value = -8196.6
$ws.Range("N43").Value = -8043.75
$ws.Range("H58").Value = 3744
$ws.Range("J58").Value = 24999
$ws.Range("L58").Value = 74997
$ws.Range("N58").Value = -75297
$ws.Range("H87").Value = 97352.5
$ws.Range("J87").Value = 97352.5
$ws.Range("L87").Value = 97352.5
$ws.Range("N87").Value = -99848.5
$ws.Range("H90").Value = 97352.5
$ws.Range("J90").Value = 97352.5
$ws.Range("L90").Value = 292057.5
$ws.Range("N90").Value = -304537.5
$ws.Range("H98").Value = 41668820
$ws.Range("I98").Value = 52633440
$ws.Range("K98").Value = 52633440
$ws.Range("M98").Value = -52631942
$ws.Range("H100").Value = 5593.4165
$ws.Range("I100").Value = 1534.4
$ws.Range("J100").Value = 25888.5
$ws.Range("K100").Value = 1534.4
$ws.Range("L100").Value = 25888.5
$ws.Range("M100").Value = -993.4000000000001
$ws.Range("N100").Value = -26970.5
$ws.Range("H112").Value = 2133.8667
$ws.Range("J112").Value = 2207.7856
$ws.Range("L112").Value = 6623.3568
$ws.Range("N112").Value = -8839.356800000001
$ws.Range("H122").Value = 41668820
$ws.Range("I122").Value = 52633440
$ws.Range("K122").Value = 157900320
$ws.Range("M122").Value = -157897870
$ws.Range("H132").Value = 1828.9286
$ws.Range("I132").Value = 1109.3182
$ws.Range("K132").Value = 3327.9546
$ws.Range("M132").Value = -797.9546
$ws.Range("H138").Value = 1907.8474
$ws.Range("J138").Value = 2111.7659
$ws.Range("L138").Value = 6335.297699999999
$ws.Range("N138").Value = -16615.2977

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 23811906
$ws.Range("I45").Value = 33335028
$ws.Range("J45").Value = 4099.8335
$ws.Range("K45").Value = 33335028
$ws.Range("L45").Value = 4099.8335
$ws.Range("M45").Value = -33334651
$ws.Range("N45").Value = -4853.8335
$ws.Range("H61").Value = 30006658
$ws.Range("I61").Value = 29417090
$ws.Range("K61").Value = 29417090
$ws.Range("M61").Value = -29416878
$ws.Range("H122").Value = 3707.2144
$ws.Range("I122").Value = 2414.7144
$ws.Range("K122").Value = 7244.1432
$ws.Range("M122").Value = -4794.1432
$ws.Range("H132").Value = 3333.6072
$ws.Range("I132").Value = 3308.963
$ws.Range("K132").Value = 9926.889000000001
$ws.Range("M132").Value = -7396.889000000001
$ws.Range("H136").Value = 30006658
$ws.Range("I136").Value = 29417090
$ws.Range("K136").Value = 88251270
$ws.Range("M136").Value = -88248720

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 25987.385
$ws.Range("J88").Value = 25987.385
$ws.Range("L88").Value = 25987.385
$ws.Range("N88").Value = -26799.385
$ws.Range("H91").Value = 25987.385
$ws.Range("J91").Value = 25987.385
$ws.Range("L91").Value = 25987.385
$ws.Range("N91").Value = -28795.385
$ws.Range("H99").Value = 2390.457
$ws.Range("I99").Value = 1875.3334
$ws.Range("K99").Value = 1875.3334
$ws.Range("M99").Value = -377.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 263.64285
$ws.Range("J7").Value = 299.42856
$ws.Range("L7").Value = 299.42856
$ws.Range("N7").Value = -525.4285600000001
$ws.Range("H16").Value = 632.4286
$ws.Range("I16").Value = 588.8461
$ws.Range("K16").Value = 588.8461
$ws.Range("M16").Value = -301.8461
$ws.Range("H58").Value = 4316.8335
$ws.Range("I58").Value = 3387.6875
$ws.Range("J58").Value = 11750
$ws.Range("K58").Value = 3387.6875
$ws.Range("L58").Value = 11750
$ws.Range("M58").Value = -3184.6875
$ws.Range("N58").Value = -12156
$ws.Range("H113").Value = 632.4286
$ws.Range("I113").Value = 588.8461
$ws.Range("K113").Value = 588.8461
$ws.Range("M113").Value = 1581.1539
$ws.Range("H132").Value = 3816.44
$ws.Range("I132").Value = 2837
$ws.Range("J132").Value = 5557.6665
$ws.Range("K132").Value = 8511
$ws.Range("L132").Value = 16672.9995
$ws.Range("M132").Value = -5981
$ws.Range("N132").Value = -21732.9995
$ws.Range("H134").Value = 3498.9285
$ws.Range("I134").Value = 3044.0908
$ws.Range("K134").Value = 9132.2724
$ws.Range("M134").Value = -6597.2724
$ws.Range("H136").Value = 4316.8335
$ws.Range("I136").Value = 3387.6875
$ws.Range("J136").Value = 11750
$ws.Range("K136").Value = 10163.0625
$ws.Range("L136").Value = 35250
$ws.Range("M136").Value = -7613.0625
$ws.Range("N136").Value = -40350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1381
$ws.Range("J113").Value = 1433.909
$ws.Range("L113").Value = 4301.727000000001
$ws.Range("N113").Value = -8641.727000000001
$ws.Range("H128").Value = 499999.34
$ws.Range("I128").Value = 499999.34
$ws.Range("K128").Value = 1499998.02
$ws.Range("M128").Value = -1495018.02
$ws.Range("H134").Value = 9203.666999999999
$ws.Range("J134").Value = 12776.88
$ws.Range("L134").Value = 38330.64
$ws.Range("N134").Value = -48470.64
$ws.Range("H140").Value = 278017.62
$ws.Range("J140").Value = 8000
$ws.Range("L140").Value = 24000
$ws.Range("N140").Value = -34360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1603.762
$ws.Range("I97").Value = 1275.2
$ws.Range("J97").Value = 2425.1667
$ws.Range("K97").Value = 1275.2
$ws.Range("L97").Value = 2425.1667
$ws.Range("M97").Value = -779.2
$ws.Range("N97").Value = -3417.1667
$ws.Range("H113").Value = 4614.846
$ws.Range("I113").Value = 5001
$ws.Range("K113").Value = 5001
$ws.Range("M113").Value = -2831
$ws.Range("H122").Value = 1561.05
$ws.Range("I122").Value = 1551.5
$ws.Range("K122").Value = 4654.5
$ws.Range("M122").Value = -2204.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4893.6
$ws.Range("I40").Value = 4490.2
$ws.Range("K40").Value = 4490.2
$ws.Range("M40").Value = -4354.2
$ws.Range("H93").Value = 40001012
$ws.Range("I93").Value = 71429340
$ws.Range("J93").Value = 1316.1818
$ws.Range("K93").Value = 71429340
$ws.Range("L93").Value = 1316.1818
$ws.Range("M93").Value = -71428092
$ws.Range("N93").Value = -3812.1818
$ws.Range("H122").Value = 6208.893
$ws.Range("I122").Value = 6038.7334
$ws.Range("J122").Value = 6405.231
$ws.Range("K122").Value = 18116.2002
$ws.Range("L122").Value = 19215.693
$ws.Range("M122").Value = -15666.2002
$ws.Range("N122").Value = -24115.693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 26079
$ws.Range("J54").Value = 31942.2
$ws.Range("L54").Value = 31942.2
$ws.Range("N54").Value = -32982.2
$ws.Range("H81").Value = 6580.2
$ws.Range("J81").Value = 49999
$ws.Range("L81").Value = 99998
$ws.Range("N81").Value = -102120
$ws.Range("H84").Value = 6580.2
$ws.Range("J84").Value = 49999
$ws.Range("L84").Value = 499990
$ws.Range("N84").Value = -510598

Write-Output "Applied all Behemoth_Profits updates"